$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 23:03"

$ws.Range("B4").Value = 1208036
$ws.Range("C4").Value = 19914
$ws.Range("D4").Value = 182490
$ws.Range("E4").Value = 956162
$ws.Range("G4").Value = 786
$ws.Range("H4").Value = 69384

$ws.Range("B9").Value = 165914
$ws.Range("C9").Value = 250
$ws.Range("E9").Value = 26279
$ws.Range("G9").Value = 69
$ws.Range("H9").Value = 6935

$ws.Range("D18").Value = 12847
$ws.Range("E18").Value = 32024

$ws.Range("D21").Value = 25200
$ws.Range("E21").Value = 2997

$ws.Range("F26").Value = 93

$ws.Range("A86").Value = "Costa de Marfil"
$ws.Range("B86").Value = 1432
$ws.Range("C86").Value = 34
$ws.Range("D86").Value = 693
$ws.Range("E86").Value = 722
$ws.Range("F86").Value = 0
$ws.Range("H86").Value = 17

$ws.Range("A87").Value = "Lituania"
$ws.Range("B87").Value = 1419
$ws.Range("C87").Value = 9
$ws.Range("D87").Value = 638
$ws.Range("E87").Value = 735
$ws.Range("F87").Value = 17
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 46

$ws.Range("A88").Value = "Eslovaquia"
$ws.Range("B88").Value = 1413
$ws.Range("C88").Value = 5
$ws.Range("D88").Value = 643
$ws.Range("E88").Value = 745
$ws.Range("F88").Value = 7
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 25

$ws.Range("A99").Value = "Niger"
$ws.Range("B99").Value = 755
$ws.Range("C99").Value = 5
$ws.Range("D99").Value = 534
$ws.Range("E99").Value = 184
$ws.Range("F99").Value = 0
$ws.Range("H99").Value = 37

$ws.Range("A100").Value = "Sri Lanka"
$ws.Range("B100").Value = 751
$ws.Range("C100").Value = 33
$ws.Range("D100").Value = 194
$ws.Range("E100").Value = 549
$ws.Range("F100").Value = 1
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 8

$ws.Range("A101").Value = "Principado de Andorra"
$ws.Range("C101").Value = 2
$ws.Range("D101").Value = 499
$ws.Range("E101").Value = 206
$ws.Range("F101").Value = 16
$ws.Range("H101").Value = 45

$ws.Range("B109").Value = 672
$ws.Range("C109").Value = 10
$ws.Range("D109").Value = 545
$ws.Range("E109").Value = 81
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 46

$ws.Range("A124").Value = "Gabon"
$ws.Range("B124").Value = 367
$ws.Range("C124").Value = 32
$ws.Range("D124").Value = 93
$ws.Range("E124").Value = 268
$ws.Range("F124").Value = 1
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 6

$ws.Range("A125").Value = "Estado de Palestina"
$ws.Range("B125").Value = 362
$ws.Range("C125").Value = 9
$ws.Range("D125").Value = 102
$ws.Range("E125").Value = 258
$ws.Range("F125").Value = 0
$ws.Range("H125").Value = 2

$ws.Range("A126").Value = "Venezuela"
$ws.Range("B126").Value = 357
$ws.Range("D126").Value = 158
$ws.Range("E126").Value = 189
$ws.Range("H126").Value = 10

$ws.Range("D153").Value = 99
$ws.Range("E153").Value = 9

$ws.Range("A155").Value = "Haiti"
$ws.Range("C155").Value = 15
$ws.Range("D155").Value = 10
$ws.Range("E155").Value = 79
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 3
$ws.Range("H155").Value = 11

$ws.Range("A156").Value = "Aruba"
$ws.Range("B156").Value = 100
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 81
$ws.Range("E156").Value = 17
$ws.Range("F156").Value = 4

$ws.Range("A157").Value = "Benin"
$ws.Range("B157").Value = 96
$ws.Range("C157").Value = 6
$ws.Range("D157").Value = 50
$ws.Range("E157").Value = 44
$ws.Range("F157").Value = 0
$ws.Range("H157").Value = 2

$ws.Range("A158").Value = "Monaco"
$ws.Range("B158").Value = 95
$ws.Range("D158").Value = 78
$ws.Range("E158").Value = 13
$ws.Range("F158").Value = 1
$ws.Range("H158").Value = 4

$ws.Range("A159").Value = "Uganda"
$ws.Range("B159").Value = 89
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 55
$ws.Range("E159").Value = 34
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 0

$ws.Range("D161").Value = 25
$ws.Range("E161").Value = 47

$ws.Range("A167").Value = "Islas Caimanes"
$ws.Range("C167").Value = 1
$ws.Range("D167").Value = 14
$ws.Range("E167").Value = 60
$ws.Range("F167").Value = 3
$ws.Range("H167").Value = 1

$ws.Range("A168").Value = "Nepal"
$ws.Range("B168").Value = 75
$ws.Range("D168").Value = 16
$ws.Range("E168").Value = 59
$ws.Range("F168").Value = 0
$ws.Range("H168").Value = 0

$ws.Range("A189").Value = "Santa Lucia"
$ws.Range("D189").Value = 15
$ws.Range("F189").Value = 0
$ws.Range("H189").Value = 0

$ws.Range("A190").Value = "Belice"
$ws.Range("D190").Value = 13
$ws.Range("F190").Value = 1
$ws.Range("H190").Value = 2

$ws.Range("A198").Value = "San Cristobal y Nieves"
$ws.Range("D198").Value = 8
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Burundi"
$ws.Range("D199").Value = 7
$ws.Range("H199").Value = 1
